{"js": "// Update the division-problem table cells in place (text-only edit; keeps\n// each run's existing formatting). Cells are addressed positionally\n// (row, column) rather than by searching for their old text, because some\n// of the new values collide with other cells' old values (e.g. a cell\n// changes from \"44\u00f74=\" to \"20\u00f79=\" while a different cell changes FROM\n// \"20\u00f79=\" to \"27\u00f77=\"), which would make a naive global find/replace\n// ambiguous or self-colliding.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, columnIndex, expectedOldText, newText]\nconst updates = [\n  [0, 0, \"30\u00f76=\", \"33\u00f77=\"],\n  [0, 1, \"44\u00f77=\", \"19\u00f77=\"],\n  [0, 2, \"27\u00f79=\", \"28\u00f77=\"],\n  [0, 3, \"66\u00f73=\", \"69\u00f72=\"],\n  [0, 4, \"10\u00f72=\", \"46\u00f73=\"],\n  [4, 0, \"27\u00f77=\", \"72\u00f76=\"],\n  [4, 1, \"44\u00f74=\", \"20\u00f79=\"],\n  [4, 2, \"50\u00f78=\", \"64\u00f78=\"],\n  [4, 3, \"58\u00f73=\", \"89\u00f73=\"],\n  [4, 4, \"47\u00f79=\", \"80\u00f76=\"],\n  [8, 0, \"41\u00f72=\", \"56\u00f72=\"],\n  [8, 1, \"11\u00f75=\", \"88\u00f79=\"],\n  [8, 2, \"15\u00f79=\", \"26\u00f79=\"],\n  [8, 3, \"65\u00f79=\", \"19\u00f79=\"],\n  [8, 4, \"22\u00f78=\", \"87\u00f75=\"],\n  [12, 0, \"46\u00f74=\", \"90\u00f78=\"],\n  [12, 1, \"80\u00f73=\", \"87\u00f73=\"],\n  [12, 2, \"45\u00f74=\", \"75\u00f78=\"],\n  [12, 3, \"17\u00f74=\", \"38\u00f79=\"],\n  [12, 4, \"23\u00f73=\", \"62\u00f74=\"],\n  [16, 0, \"35\u00f73=\", \"29\u00f78=\"],\n  [16, 1, \"20\u00f79=\", \"27\u00f77=\"],\n  [16, 2, \"71\u00f72=\", \"25\u00f79=\"],\n  [16, 3, \"14\u00f74=\", \"31\u00f76=\"],\n  [16, 4, \"96\u00f74=\", \"57\u00f73=\"],\n];\n\nconst cells = updates.map(([r, c]) => table.getCell(r, c));\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\nfor (let i = 0; i < updates.length; i++) {\n  const [, , oldText, newText] = updates[i];\n  const cell = cells[i];\n  if (cell.body.text !== oldText) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: expected \"${oldText}\" but found \"${cell.body.text}\"`\n    );\n  }\n  cell.body.getRange().insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem table cells in place (text-only edit; keeps\n# each cell's existing run formatting). Cells are addressed positionally\n# (1-based Table.Cell(row, column)) rather than by searching for their old\n# text, because some of the new values collide with other cells' old\n# values (e.g. a cell changes from \"44\u00f74=\" to \"20\u00f79=\" while a different\n# cell changes FROM \"20\u00f79=\" to \"27\u00f77=\"), which would make a naive\n# Find/Replace across the whole document ambiguous or self-colliding.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each entry: row, column (1-based), expected old text, new text\n$updates = @(\n    @(1, 1, \"30\u00f76=\", \"33\u00f77=\"),\n    @(1, 2, \"44\u00f77=\", \"19\u00f77=\"),\n    @(1, 3, \"27\u00f79=\", \"28\u00f77=\"),\n    @(1, 4, \"66\u00f73=\", \"69\u00f72=\"),\n    @(1, 5, \"10\u00f72=\", \"46\u00f73=\"),\n    @(5, 1, \"27\u00f77=\", \"72\u00f76=\"),\n    @(5, 2, \"44\u00f74=\", \"20\u00f79=\"),\n    @(5, 3, \"50\u00f78=\", \"64\u00f78=\"),\n    @(5, 4, \"58\u00f73=\", \"89\u00f73=\"),\n    @(5, 5, \"47\u00f79=\", \"80\u00f76=\"),\n    @(9, 1, \"41\u00f72=\", \"56\u00f72=\"),\n    @(9, 2, \"11\u00f75=\", \"88\u00f79=\"),\n    @(9, 3, \"15\u00f79=\", \"26\u00f79=\"),\n    @(9, 4, \"65\u00f79=\", \"19\u00f79=\"),\n    @(9, 5, \"22\u00f78=\", \"87\u00f75=\"),\n    @(13, 1, \"46\u00f74=\", \"90\u00f78=\"),\n    @(13, 2, \"80\u00f73=\", \"87\u00f73=\"),\n    @(13, 3, \"45\u00f74=\", \"75\u00f78=\"),\n    @(13, 4, \"17\u00f74=\", \"38\u00f79=\"),\n    @(13, 5, \"23\u00f73=\", \"62\u00f74=\"),\n    @(17, 1, \"35\u00f73=\", \"29\u00f78=\"),\n    @(17, 2, \"20\u00f79=\", \"27\u00f77=\"),\n    @(17, 3, \"71\u00f72=\", \"25\u00f79=\"),\n    @(17, 4, \"14\u00f74=\", \"31\u00f76=\"),\n    @(17, 5, \"96\u00f74=\", \"57\u00f73=\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $oldText = $u[2]\n    $newText = $u[3]\n\n    $cell = $t.Cell($row, $col)\n    $range = $cell.Range\n    $current = $range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -ne $oldText) {\n        throw \"Unexpected text at row $row, col $col`: expected '$oldText' but found '$current'\"\n    }\n\n    $range.Text = $newText\n}\n"}
